$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28; existing rows 28..60 shift down to 29..61
$ws.Rows.Item(28).Insert()

# Fill in the new row 28 with the data
$ws.Cells.Item(28, 1).Value = 3
$ws.Cells.Item(28, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44763
$ws.Cells.Item(28, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(28, 5).Value = 5
$ws.Cells.Item(28, 6).Value = 100112035
$ws.Cells.Item(28, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 80
$ws.Cells.Item(28, 11).Value = 14000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 14500
$ws.Cells.Item(28, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(28, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(28, 16).Value = 967
$ws.Cells.Item(28, 17).Value = 15
$ws.Cells.Item(28, 18).Value = "Hortaliza"
